# Add SQF-RNN results and evaluation
# This script reproduces, on "S-score-0" / "S-score-2" / "S-score-3", the
# insertion of 6 new result rows (sqf-rnn with TP / sqf-rnn with TP, TCC /
# nnqf with TP / nnqf with TP, TCC / nnqf with TCC) between the existing
# "sqf-rnn" / "nnqf" and "deepar" rows, plus a new labeled "task" row.

$wb = $excel.ActiveWorkbook

$ws0 = $wb.Worksheets.Item("S-score-0")
$ws2 = $wb.Worksheets.Item("S-score-2")
$ws3 = $wb.Worksheets.Item("S-score-3")

$cols = @("B","C","D","E","F","G","H","I","J","K","L","M")

# ---------------------------------------------------------------------
# 1) Insert the new rows. Order matters: insert into the most-dependent
#    sheet first (S-score-3 depends on S-score-2 depends on S-score-0)
#    so that each insert's automatic reference-shifting of formulas that
#    already point at the lower sheet keeps everything aligned by the
#    time we reach S-score-0.
# ---------------------------------------------------------------------
foreach ($ws in @($ws3, $ws2, $ws0)) {
    $ws.Rows("25:26").Insert()
    $ws.Rows("28:30").Insert()
}

# ---------------------------------------------------------------------
# 2) Prime the shared-string table in the same order the new labels were
#    first typed (so new entries land at the same indices as the target
#    workbook: 339 nnqf with TP, 340 sqf-rnn with TP TCC, 341 sqf-rnn
#    with TP, 342 nnqf with TP TCC, 343 nnqf with TCC, 344 task).
# ---------------------------------------------------------------------
$ws0.Range("A28").Value = "nnqf with TP"
$ws0.Range("A26").Value = "sqf-rnn with TP, TCC"
$ws0.Range("A25").Value = "sqf-rnn with TP"
$ws0.Range("A29").Value = "nnqf with TP, TCC"
$ws0.Range("A30").Value = "nnqf with TCC"
$ws0.Range("A33").Value = "task"

# ---------------------------------------------------------------------
# 3) S-score-0 : raw data sheet
# ---------------------------------------------------------------------

# Row 24 "sqf-rnn" : drop the stray note that used to live in O24
$ws0.Range("O24").ClearContents()

# Row 25 "sqf-rnn with TP"
$ws0.Range("L25").Value = 0.01665

# Row 26 "sqf-rnn with TP, TCC"
$ws0.Range("B26").Value = 0.02313
$ws0.Range("C26").Value = 0.02921
$ws0.Range("D26").Value = 0.02491
$ws0.Range("E26").Value = 0.02445
$ws0.Range("F26").Value = 0.02055
$ws0.Range("G26").Value = 0.02057
$ws0.Range("H26").Value = 0.01925
$ws0.Range("I26").Value = 0.01942
$ws0.Range("J26").Value = 0.02165
$ws0.Range("K26").Value = 0.02231
$ws0.Range("L26").Value = 0.01207
$ws0.Range("M26").Value = 0.01916
$ws0.Range("O26").Value = "ensemble_count=7, epochs=7"

# Row 27 "nnqf" (previously row 25, now shifted down with its own new siblings)
$ws0.Range("A27").Value = "nnqf"
$ws0.Range("F27").Value = 0.03348

# Row 28 "nnqf with TP"
$ws0.Range("F28").Value = 0.03376
$ws0.Range("H28").Value = 0.02951
$ws0.Range("L28").Value = 0.0219

# Row 29 "nnqf with TP, TCC"
$ws0.Range("F29").Value = 0.03363
$ws0.Range("L29").Value = 0.0228

# Row 30 "nnqf with TCC"
$ws0.Range("F30").Value = 0.03318
$ws0.Range("L30").Value = 0.02279

# Row 31 "deepar" (previously row 26) - label only, no data change
$ws0.Range("A31").Value = "deepar"

# Row 33 "task" header row (previously row 28, now labeled in column A)
$ws0.Range("B33").Value = 4
$ws0.Range("C33").Value = 5
$ws0.Range("D33").Value = 6
$ws0.Range("E33").Value = 7
$ws0.Range("F33").Value = 8
$ws0.Range("G33").Value = 9
$ws0.Range("H33").Value = 10
$ws0.Range("I33").Value = 11
$ws0.Range("J33").Value = 12
$ws0.Range("K33").Value = 13
$ws0.Range("L33").Value = 14
$ws0.Range("M33").Value = 15

# ---------------------------------------------------------------------
# 4) S-score-2 : mirrors S-score-0 via ='S-score-0'!X{row} formulas
# ---------------------------------------------------------------------
$ws2.Range("A25").Value = "sqf-rnn with TP"
$ws2.Range("A26").Value = "sqf-rnn with TP, TCC"
$ws2.Range("A27").Value = "nnqf"
$ws2.Range("A28").Value = "nnqf with TP"
$ws2.Range("A29").Value = "nnqf with TP, TCC"
$ws2.Range("A30").Value = "nnqf with TCC"
$ws2.Range("A31").Value = "deepar"

for ($r = 25; $r -le 31; $r++) {
    foreach ($c in $cols) {
        $ws2.Range("$c$r").Formula = "='S-score-0'!$c$r"
    }
}

# ---------------------------------------------------------------------
# 5) S-score-3 : mirrors S-score-2 via 1-'S-score-2'!X{row}/'S-score-2'!X$4
#    plus the weighted-average column N and the (new) RANK column O.
# ---------------------------------------------------------------------
$ws3.Range("A25").Value = "sqf-rnn with TP"
$ws3.Range("A26").Value = "sqf-rnn with TP, TCC"
$ws3.Range("A27").Value = "nnqf"
$ws3.Range("A28").Value = "nnqf with TP"
$ws3.Range("A29").Value = "nnqf with TP, TCC"
$ws3.Range("A30").Value = "nnqf with TCC"
$ws3.Range("A31").Value = "deepar"

for ($r = 25; $r -le 31; $r++) {
    foreach ($c in $cols) {
        $ws3.Range("$c$r").Formula = "=1-'S-score-2'!$c$r/'S-score-2'!$c`$4"
    }
    $ws3.Range("N$r").Formula = "=SUMPRODUCT(B`$1:M`$1,B${r}:M$r)/SUM(B`$1:M`$1)"
}

# New RANK column for rows 24-31 (row 24 "sqf-rnn" previously had no O cell)
for ($r = 24; $r -le 31; $r++) {
    $ws3.Range("O$r").Formula = "=RANK(N$r,(N`$2:N`$22,N$r))"
}

# ---------------------------------------------------------------------
# 6) Restore the view/selection state shown in the target workbook.
# ---------------------------------------------------------------------
$ws3.Activate()
$ws3.Range("M21").Select()

$ws2.Activate()
$ws2.Range("A24:A31").Select()

$ws0.Activate()
$ws0.Range("O26").Select()
